$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("E5").Value = 1300
$ws.Range("H5").Value = 45251.04190482639
$ws.Range("J5").Value = '10/28/23 17:20'
$ws.Range("K5").Value = '10/28/23 17:20'
$ws.Range("M5").Value = '$1,300 as of 10/28/2023 3:20:32 PM'
$ws.Range("N5").Value = 1500

# Row 6
$ws.Range("A6").Value = 'L662336'
$ws.Range("C6").Value = 'SB#4 MONA MARKET'
$ws.Range("E6").Value = 1440
$ws.Range("H6").Value = 45234.04190482639
$ws.Range("I6").ClearContents()
$ws.Range("J6").Value = '10/29/23 14:38'
$ws.Range("K6").Value = '10/29/23 14:38'
$ws.Range("L6").Value = 120
$ws.Range("M6").Value = '$1,480 as of 10/29/2023 9:25:15 AM'
$ws.Range("N6").Value = 1480

# Row 7
$ws.Range("A7").Value = 'L647934'
$ws.Range("C7").Value = 'SB #6'
$ws.Range("E7").Value = 1940
$ws.Range("J7").Value = '04/06/23 22:10'
$ws.Range("K7").Value = '04/06/23 22:05'
$ws.Range("L7").Value = 20
$ws.Range("M7").Value = '$1,940 as of 4/6/2023 8:05:45 PM'
$ws.Range("N7").Value = 1960

# Row 8
$ws.Range("A8").Value = 'L474792'
$ws.Range("C8").Value = 'NICK SHELL SERVICE'
$ws.Range("E8").Value = 2140
$ws.Range("H8").Value = 45243.04190482639
$ws.Range("I8").ClearContents()
$ws.Range("J8").Value = '10/29/23 09:18'
$ws.Range("K8").Value = '10/29/23 09:18'
$ws.Range("M8").Value = '$2,140 as of 10/29/2023 7:18:35 AM'
$ws.Range("N8").Value = 2160

# Row 9
$ws.Range("A9").Value = 'L476340'
$ws.Range("C9").Value = 'DONUT & SANDWICH'
$ws.Range("E9").Value = 2360
$ws.Range("H9").Value = 45236.04190482639
$ws.Range("J9").Value = '10/29/23 12:47'
$ws.Range("K9").Value = '10/29/23 12:47'
$ws.Range("M9").Value = '$2,360 as of 10/29/2023 10:47:34 AM'
$ws.Range("N9").Value = 2460

# Row 10
$ws.Range("A10").Value = 'L678988'
$ws.Range("C10").Value = 'PAYELESS MARKET'
$ws.Range("E10").Value = 2400
$ws.Range("H10").ClearContents()
$ws.Range("I10").Value = 'ATM Inactive greater than 2000 minutes'
$ws.Range("J10").Value = '07/20/23 20:09'
$ws.Range("K10").Value = '07/20/23 20:09'
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = '$2,400 as of 7/20/2023 6:09:40 PM'
$ws.Range("N10").Value = 2500

# Row 11
$ws.Range("A11").Value = 'LK561655'
$ws.Range("C11").Value = 'CRENSHAW CRAVOR #2'
$ws.Range("E11").Value = 2780
$ws.Range("H11").ClearContents()
$ws.Range("I11").Value = 'ATM Inactive greater than 48 minutes'
$ws.Range("J11").Value = '01/23/20 08:24'
$ws.Range("K11").Value = '01/23/20 08:24'
$ws.Range("M11").Value = '$2,780 as of 1/23/2020 6:24:32 AM'
$ws.Range("N11").Value = 2800

# Row 12
$ws.Range("A12").Value = 'L488595'
$ws.Range("C12").Value = 'N S MART'
$ws.Range("E12").Value = 3140
$ws.Range("H12").Value = 45307.04190482639
$ws.Range("I12").ClearContents()
$ws.Range("J12").Value = '10/28/23 21:40'
$ws.Range("K12").Value = '10/28/23 21:40'
$ws.Range("M12").Value = '$3,140 as of 10/28/2023 7:40:28 PM'
$ws.Range("N12").Value = 3220

# Row 13
$ws.Range("A13").Value = 'LK236828'
$ws.Range("C13").Value = 'WORLDWIDE AUTOMOTIVE'
$ws.Range("E13").Value = 3560
$ws.Range("H13").Value = 45257.04190482639
$ws.Range("I13").ClearContents()
$ws.Range("J13").Value = '10/29/23 16:35'
$ws.Range("K13").Value = '10/29/23 16:35'
$ws.Range("L13").Value = 20
$ws.Range("M13").Value = '$3,620 as of 10/28/2023 10:37:24 AM'
$ws.Range("N13").Value = 3620

# Row 14
$ws.Range("A14").Value = 'L697590'
$ws.Range("C14").Value = 'S B MARKET ST'
$ws.Range("E14").Value = 3720
$ws.Range("H14").Value = 45265.04190482639
$ws.Range("J14").Value = '10/29/23 14:07'
$ws.Range("K14").Value = '10/29/23 14:07'
$ws.Range("M14").Value = '$3,740 as of 10/28/2023 2:56:53 PM'
$ws.Range("N14").Value = 3740

# Row 15
$ws.Range("A15").Value = 'L475182'
$ws.Range("C15").Value = 'LA ESQUINA DE ORO'
$ws.Range("E15").Value = 3800
$ws.Range("H15").ClearContents()
$ws.Range("I15").Value = 'ATM Inactive greater than 48 minutes'
$ws.Range("J15").Value = '09/16/20 16:57'
$ws.Range("K15").Value = '09/15/20 23:38'
$ws.Range("M15").Value = '$3,800 as of 9/16/2020 1:28:00 PM'
$ws.Range("N15").Value = 3800

# Row 16
$ws.Range("A16").Value = 'L688966'
$ws.Range("C16").Value = 'S B WESTERN 108TH MARKET'
$ws.Range("E16").Value = 3980
$ws.Range("H16").Value = 45272.04190482639
$ws.Range("J16").Value = '10/29/23 15:09'
$ws.Range("K16").Value = '10/29/23 15:09'
$ws.Range("L16").Value = 80
$ws.Range("M16").Value = '$4,100 as of 10/29/2023 11:15:25 AM'
$ws.Range("N16").Value = 4080

# Row 17
$ws.Range("A17").Value = 'L688961'
$ws.Range("C17").Value = 'MONA MART'
$ws.Range("E17").Value = 4000
$ws.Range("H17").ClearContents()
$ws.Range("I17").Value = 'ATM Inactive greater than 2000 minutes'
$ws.Range("J17").Value = '10/17/23 13:26'
$ws.Range("K17").Value = '10/17/23 13:00'
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = '$4,000 as of 10/17/2023 11:00:09 AM'
$ws.Range("N17").Value = 4000

# Row 18
$ws.Range("A18").Value = 'L475090'
$ws.Range("C18").Value = 'S.B. 2'
$ws.Range("E18").Value = 5000
$ws.Range("H18").Value = 45238.04190482639
$ws.Range("J18").Value = '10/29/23 18:00'
$ws.Range("K18").Value = '10/29/23 18:00'
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = '$5,080 as of 10/28/2023 4:33:08 PM'
$ws.Range("N18").Value = 5020

# Row 19
$ws.Range("A19").Value = 'L474746'
$ws.Range("C19").Value = 'ZACATES MARKET'
$ws.Range("E19").Value = 5440
$ws.Range("H19").Value = 45255.04190482639
$ws.Range("J19").Value = '10/29/23 17:21'
$ws.Range("K19").Value = '10/29/23 17:21'
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = '$5,620 as of 10/29/2023 11:14:47 AM'
$ws.Range("N19").Value = 5620

# Row 20
$ws.Range("A20").Value = 'LK864765'
$ws.Range("C20").Value = 'SKY LIQUOR'
$ws.Range("E20").Value = 5520
$ws.Range("H20").Value = 45236.04190482639
$ws.Range("J20").Value = '10/29/23 12:55'
$ws.Range("K20").Value = '10/29/23 02:29'
$ws.Range("L20").Value = 80
$ws.Range("M20").Value = '$5,600 as of 10/29/2023 10:00:03 AM'
$ws.Range("N20").Value = 5600

# Row 21
$ws.Range("A21").Value = 'L704741'
$ws.Range("C21").Value = 'W ADAMS COIN LAUNDRY'
$ws.Range("E21").Value = 5800
$ws.Range("H21").Value = 45237.04190482639
$ws.Range("J21").Value = '10/29/23 17:36'
$ws.Range("K21").Value = '10/29/23 17:36'
$ws.Range("L21").Value = 40
$ws.Range("M21").Value = '$5,860 as of 10/29/2023 9:54:54 AM'
$ws.Range("N21").Value = 5840

# Row 22
$ws.Range("A22").Value = 'L474817'
$ws.Range("C22").Value = 'SAFETY MARKET'
$ws.Range("E22").Value = 6700
$ws.Range("H22").Value = 45244.04190482639
$ws.Range("I22").ClearContents()
$ws.Range("J22").Value = '10/29/23 18:10'
$ws.Range("K22").Value = '10/29/23 18:10'
$ws.Range("L22").Value = 120
$ws.Range("M22").Value = '$7,280 as of 10/29/2023 8:50:55 AM'
$ws.Range("N22").Value = 6800

# Row 23
$ws.Range("E23").Value = 6860
$ws.Range("H23").Value = 45268.04190482639
$ws.Range("J23").Value = '10/29/23 12:02'
$ws.Range("K23").Value = '10/29/23 12:02'
$ws.Range("M23").Value = '$6,860 as of 10/29/2023 10:02:36 AM'
$ws.Range("N23").Value = 6920

# Row 24
$ws.Range("A24").Value = 'LK891176'
$ws.Range("C24").Value = '98 DISCOUNT STORE'
$ws.Range("E24").Value = 6920
$ws.Range("H24").Value = 45236.04190482639
$ws.Range("J24").Value = '10/29/23 17:28'
$ws.Range("K24").Value = '10/29/23 17:28'
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = '$7,280 as of 10/29/2023 10:20:08 AM'
$ws.Range("N24").Value = 6960

# Row 25
$ws.Range("A25").Value = 'L682801'
$ws.Range("C25").Value = 'SB#5'
$ws.Range("E25").Value = 7840
$ws.Range("H25").ClearContents()
$ws.Range("I25").Value = 'ATM Inactive greater than 2000 minutes'
$ws.Range("J25").Value = '09/28/23 15:22'
$ws.Range("K25").Value = '09/28/23 12:14'
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = '$7,840 as of 9/28/2023 12:31:50 PM'
$ws.Range("N25").Value = 7840

# Row 26
$ws.Range("A26").Value = 'L697589'
$ws.Range("C26").Value = 'S B DISCOUNT MART'
$ws.Range("E26").Value = 13240
$ws.Range("H26").Value = 45242.04190482639
$ws.Range("J26").Value = '10/29/23 18:14'
$ws.Range("K26").Value = '10/29/23 17:49'
$ws.Range("L26").Value = 60
$ws.Range("M26").Value = '$13,380 as of 10/29/2023 10:17:18 AM'
$ws.Range("N26").Value = 13240

# Row 27
$ws.Range("E27").Value = 14760
$ws.Range("H27").Value = 45255.04190482639
$ws.Range("J27").Value = '10/29/23 15:29'
$ws.Range("K27").Value = '10/29/23 15:29'
$ws.Range("M27").Value = '$14,780 as of 10/29/2023 10:00:13 AM'
$ws.Range("N27").Value = 14780

# Row 28
$ws.Range("E28").Value = 114640
